# Automatische test-sync: 2025-06-30 20:19:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 19

$logs.Cells.Item($row, 1).Value = "Heeft u informatie over zakelijke kortingen voor wederverkopers?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #19: Heeft u informatie over zakelijke kortingen voor wederverkopers?"
$logs.Cells.Item($row, 4).Value = "Productinformatie"
$logs.Cells.Item($row, 5).Value = "Beste klant,
Bedankt voor uw interesse in onze zakelijke kortingen voor wederverkopers. Voor meer informatie over onze zakelijke kortingen en de voorwaarden kunt u contact opnemen met onze verkoopafdeling via sales@bedrijfsnaam.nl. Zij kunnen u voorzien van alle benodigde informatie en u verder helpen met uw aanvraag.
Met vriendelijke groet,
[Naam bedrijf] E-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-30 20:19:21"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"
$logs.Cells.Item($row, 10).Value = "Nee"

# Setting a multi-line value auto-expands the row height; reset it back to
# the sheet default so the new row matches the rest of the sheet.
$logs.Rows.Item($row).AutoFit()

# Extend the conditional formatting ranges so they cover the new row too.
foreach ($col in @("D","G","H","I","J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "18")
    $newRange = $logs.Range($col + "2:" + $col + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Productinformatie" (row 3, column B)
$dashboard.Cells.Item(3, 2).Value = 5
